$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add F and G values for rows 393-422
$ws.Cells.Item(393, 6).Value = 308016
$ws.Cells.Item(393, 7).Value = 1240
$ws.Cells.Item(394, 6).Value = 165775
$ws.Cells.Item(394, 7).Value = 623
$ws.Cells.Item(395, 6).Value = 752222
$ws.Cells.Item(395, 7).Value = 1950
$ws.Cells.Item(396, 6).Value = 166408
$ws.Cells.Item(396, 7).Value = 549
$ws.Cells.Item(397, 6).Value = 107795
$ws.Cells.Item(397, 7).Value = 639
$ws.Cells.Item(398, 6).Value = 298822
$ws.Cells.Item(398, 7).Value = 1468
$ws.Cells.Item(399, 6).Value = 200323
$ws.Cells.Item(399, 7).Value = 967
$ws.Cells.Item(400, 6).Value = 148992
$ws.Cells.Item(400, 7).Value = 764
$ws.Cells.Item(401, 6).Value = 272418
$ws.Cells.Item(401, 7).Value = 936
$ws.Cells.Item(402, 6).Value = 721927
$ws.Cells.Item(402, 7).Value = 1392
$ws.Cells.Item(403, 6).Value = 353822
$ws.Cells.Item(403, 7).Value = 734
$ws.Cells.Item(404, 6).Value = 224126
$ws.Cells.Item(404, 7).Value = 914
$ws.Cells.Item(405, 6).Value = 174024
$ws.Cells.Item(405, 7).Value = 693
$ws.Cells.Item(406, 6).Value = 170944
$ws.Cells.Item(406, 7).Value = 680
$ws.Cells.Item(407, 6).Value = 158075
$ws.Cells.Item(407, 7).Value = 673
$ws.Cells.Item(408, 6).Value = 304552
$ws.Cells.Item(408, 7).Value = 835
$ws.Cells.Item(409, 6).Value = 708499
$ws.Cells.Item(409, 7).Value = 1007
$ws.Cells.Item(410, 6).Value = 364399
$ws.Cells.Item(410, 7).Value = 635
$ws.Cells.Item(411, 6).Value = 225424
$ws.Cells.Item(411, 7).Value = 828
$ws.Cells.Item(412, 6).Value = 176157
$ws.Cells.Item(412, 7).Value = 646
$ws.Cells.Item(413, 6).Value = 149578
$ws.Cells.Item(413, 7).Value = 658
$ws.Cells.Item(414, 6).Value = 148839
$ws.Cells.Item(414, 7).Value = 563
$ws.Cells.Item(415, 6).Value = 307891
$ws.Cells.Item(415, 7).Value = 694
$ws.Cells.Item(416, 6).Value = 671641
$ws.Cells.Item(416, 7).Value = 931
$ws.Cells.Item(417, 6).Value = 342466
$ws.Cells.Item(417, 7).Value = 589
$ws.Cells.Item(418, 6).Value = 202140
$ws.Cells.Item(418, 7).Value = 700
$ws.Cells.Item(419, 6).Value = 149297
$ws.Cells.Item(419, 7).Value = 510
$ws.Cells.Item(420, 6).Value = 138714
$ws.Cells.Item(420, 7).Value = 500
$ws.Cells.Item(421, 6).Value = 152971
$ws.Cells.Item(421, 7).Value = 532
$ws.Cells.Item(422, 6).Value = 298331
$ws.Cells.Item(422, 7).Value = 645

# Add new row 453 with full data
$ws.Cells.Item(453, 1).Value = 44347
$ws.Cells.Item(453, 1).NumberFormat = "yyyy-mm-dd"
$ws.Cells.Item(453, 2).Value = 389866
$ws.Cells.Item(453, 3).Value = 7086
$ws.Cells.Item(453, 4).Value = 145
$ws.Cells.Item(453, 5).Value = 12353
$ws.Cells.Item(453, 6).Value = 57936
$ws.Cells.Item(453, 7).Value = 174
